# Updates cryptos list figures (prices / volume change percentages) and
# swaps the InjectiveProtocol / VeChain rows (33 & 34).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 33 / 34 swap (coin name, link, price, volume) ---
$ws.Range("B33").Value = "VeChain"
$ws.Range("C33").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0498"
$ws.Range("E33").Value = "  +11.38%  "

$ws.Range("B34").Value = "InjectiveProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "37.17"
$ws.Range("E34").Value = "  +7.27%  "

# --- Price / Volume(1h) updates for remaining rows ---

# Row 2 - Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "53.278.73"
$ws.Range("E2").Value = "  +3.59%  "

# Row 3 - Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.149.77"
$ws.Range("E3").Value = "  +3.34%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.03%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "398.23"
$ws.Range("E5").Value = "  +3.40%  "

# Row 6 - Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "107.56"
$ws.Range("E6").Value = "  +4.46%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.00%  "

# Row 9 - Cardano
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.608"
$ws.Range("E9").Value = "  +4.06%  "

# Row 10 - Avalanche
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.79"
$ws.Range("E10").Value = "  +5.58%  "

# Row 11 - TRON
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.139"
$ws.Range("E11").Value = "  +1.32%  "

# Row 12 - Dogecoin
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0869"
$ws.Range("E12").Value = "  +0.91%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.651.72"
$ws.Range("E13").Value = "  +3.32%  "

# Row 14 - Chainlink
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.00"
$ws.Range("E14").Value = "  +2.55%  "

# Row 15 - Polkadot
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.96"
$ws.Range("E15").Value = "  +2.56%  "

# Row 16 - Polygon
$ws.Range("E16").Value = "  +8.98%  "

# Row 17 - WrappedEther
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.155.31"
$ws.Range("E17").Value = "  +2.95%  "

# Row 18 - Uniswap
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.58"
$ws.Range("E18").Value = "  -0.23%  "

# Row 19 - WrappedBTC
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "53.287.12"
$ws.Range("E19").Value = "  +3.41%  "

# Row 20 - ImmutableX
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.28"
$ws.Range("E20").Value = "  +4.31%  "

# Row 21 - InternetComputer(DFINITY)
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.83"
$ws.Range("E21").Value = "  +3.44%  "

# Row 22 - ShibaInu
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0973"
$ws.Range("E22").Value = "  +0.72%  "

# Row 23 - Litecoin
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.80"
$ws.Range("E23").Value = "  +1.02%  "

# Row 24 - BitcoinCash
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "270.79"
$ws.Range("E24").Value = "  +1.12%  "

# Row 25 - PancakeSwap
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.21"
$ws.Range("E25").Value = "  +1.82%  "

# Row 26 - Filecoin
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.09"
$ws.Range("E26").Value = "  -1.42%  "

# Row 27 - EthereumClassic
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "27.71"
$ws.Range("E27").Value = "  +3.02%  "

# Row 28 - RenderToken
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.39"
$ws.Range("E28").Value = "  +2.27%  "

# Row 29 - Kaspa
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.169"
$ws.Range("E29").Value = "  -0.85%  "

# Row 30 - Dai
$ws.Range("E30").Value = "  -0.08%  "

# Row 31 - Hedera
$ws.Range("E31").Value = "  +2.62%  "

# Row 32 - Cosmos
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.97"
$ws.Range("E32").Value = "  +7.03%  "

# Row 35 - Toncoin
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.07"
$ws.Range("E35").Value = "  +0.28%  "

# Row 36 - OKB
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "50.31"
$ws.Range("E36").Value = "  -0.21%  "

# Row 37 - LidoDAOToken
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.61"
$ws.Range("E37").Value = "  +8.28%  "

# Row 38 - FirstDigitalUSD
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"

# Row 39 - Stacks
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.78"
$ws.Range("E39").Value = "  +8.93%  "

# Row 40 - NEARProtocol
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.17"
$ws.Range("E40").Value = "  +11.31%  "

# Row 41 - TheGraph
$ws.Range("E41").Value = "  -0.25%  "

# Row 42 - Celestia
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "17.37"
$ws.Range("E42").Value = "  +2.57%  "

# Row 43 - ARBITRUM
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.90"
$ws.Range("E43").Value = "  +2.08%  "

# Row 44 - Monero
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "130.46"
$ws.Range("E44").Value = "  +4.35%  "

# Row 45 - Stellar
$ws.Range("E45").Value = "  +1.30%  "

# Row 46 - EnergySwap
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.37"
$ws.Range("E46").Value = "  +2.47%  "

# Row 47 - WEMIXToken
$ws.Range("E47").Value = "  -1.14%  "

# Row 48 - Maker
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.091.10"
$ws.Range("E48").Value = "  +2.94%  "

# Row 49 - ApeXProtocol
$ws.Range("E49").Value = "  +0.00%  "

# Row 50 - FlareNetwork
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0504"
$ws.Range("E50").Value = "  +24.01%  "

# Row 51 - BEAM
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0332"
$ws.Range("E51").Value = "  +4.29%  "

